$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "3.4.101a"
$ws.Range("B5").Value = "VCASU00020"

$ws.Range("A6").Value = "3.4.101a"
$ws.Range("B6").Value = "VCASU00030"

$ws.Range("A7").Value = "3.4.101a"
$ws.Range("B7").Value = "VCASU00040"
